# Remove duplication of sample_id column in dcp2T1
# Clears the duplicated sample_id values from column B (rows 6-34)
# on the "Tier 1_obs" sheet, since column B duplicated the per-sample
# identifiers already represented elsewhere.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Tier 1_obs")

$ws.Range("B6:B34").ClearContents()
